# January-2021-RetailerWiseOrder.xlsx — "data updated till 23 Jan 9AM"
#
# 1. Update the raw daily-order input cells (columns AA:AM) that carry new
#    "Auto"/"Manual" dispatch figures. Row-2 totals and the F/G formula
#    columns on each row recalc automatically from these.
# 2. Move/retext the Vijay comment that was anchored on AE7 to AD7.
# 3. Add the new Vijay comment on AC53.
# 4. Remove the now-obsolete Vijay comment on AE72.
# 5. Update the active selection to D91 to match where the author left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Cell value updates --------------------------------------------------

$cellUpdates = @{
    "AD4"  = 1040
    "AC5"  = 1040
    "AC7"  = 3120
    "AD7"  = 6240
    "AD8"  = 2080
    "AD13" = 5200
    "AC15" = 2080
    "AD17" = 3120
    "AD19" = 2080
    "AE23" = 3120
    "AC25" = 3120
    "AC27" = 2080
    "AD30" = 5200
    "AD35" = 3120
    "AE40" = 2080
    "AD41" = 3120
    "AC46" = 3120
    "AD48" = 3120
    "AD52" = 1040
    "AC53" = 2080
    "AD53" = 1040
    "AC55" = 5200
    "AD59" = 2080
    "AE60" = 2080
    "AD63" = 2080
    "AC64" = 1040
    "AC65" = 5200
    "AE66" = 3120
    "AD69" = 1040
    "AE69" = 1040
    "AD70" = 3120
    "AD72" = 2080
    "AC75" = 2080
    "AD80" = 1040
    "AC82" = 5200
    "AD90" = 5200
    "AD91" = 5200
}

foreach ($ref in $cellUpdates.Keys) {
    $ws.Range($ref).Value = $cellUpdates[$ref]
}

# Row 90 also gains a retailer/branch label in column D (reuses the existing
# "Arwal" text already used elsewhere in column D, e.g. D82).
$ws.Range("D90").Value = "Arwal"

# --- 2. Move comment AE7 -> AD7, retexted -----------------------------------

$movedComment = $ws.Range("AE7").Comment
$movedComment.Delete()
$ws.Range("AD7").AddComment("Vijay:`n4160-Auto`n2080-Auto") | Out-Null

# --- 3. New comment on AC53 --------------------------------------------------

$ws.Range("AC53").AddComment("Vijay:`n1040-Auto`n1040-Auto") | Out-Null

# --- 4. Drop the stale AE72 comment -----------------------------------------

$ws.Range("AE72").Comment.Delete()

# --- 5. Leave the selection where the author left off ------------------------

$ws.Activate()
$ws.Range("D91").Select() | Out-Null
